$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price (column D) values are plain decimal numbers
# (e.g. "562.35"). The source workbook stores every Price cell as text, so
# without help Excel would silently reinterpret these as numeric cells.
# Temporarily force a text format on just those cells, assign the values,
# then restore their original (default/"Normal") style.
$forceTextCells = @("D5", "D6", "D12", "D14", "D20", "D22", "D24", "D26", "D27", "D31", "D37", "D39", "D40", "D41", "D45", "D47", "D48", "D49", "D50")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.069.61"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "2.423.67"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D5").Value = "562.35"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "143.71"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "2.422.65"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "5.21"
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "26.17"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").Value = "2.859.71"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "61.961.21"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "2.415.86"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "323.14"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").Value = "6.82"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "67.42"
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").Value = "8.59"
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("D27").Value = "558.27"
$ws.Range("E27").Value = "  -4.53%  "
$ws.Range("D28").Value = "2.540.19"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").Value = "0.0₃0934"
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("D31").Value = "8.20"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("E32").Value = "  -4.75%  "
$ws.Range("E33").Value = "  -2.37%  "
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "4.75"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "5.48"
$ws.Range("E39").Value = "  -4.40%  "
$ws.Range("D40").Value = "151.95"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").Value = "18.67"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").Value = "147.17"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "0.0529"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").Value = "19.96"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("D49").Value = "0.595"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "0.0923"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("E51").Value = "  -0.55%  "

foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).Style = "Normal"
}
